# Update the benchmark results table (columns B:G, rows 2:15) on the "chart" sheet
# with refreshed timing numbers, as produced by a re-run of the benchmark after the
# "optimize redundant generation when tail is empty" perf change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chart")

$data = @(
    @(4038, 4018, 300,  1865, 1806, 99),
    @(4964, 4974, 307,  1864, 1795, 102),
    @(4286, 4035, 560,  2582, 2129, 514),
    @(4171, 4030, 419,  3446, 3318, 115),
    @(5101, 4971, 424,  3396, 3319, 114),
    @(4292, 3997, 586,  4319, 4001, 417),
    @(4107, 3939, 452,  3337, 3224, 119),
    @(4177, 4038, 442,  3366, 3262, 121),
    @(4894, 4095, 1115, 4633, 3923, 753),
    @(5364, 4948, 712,  3441, 3301, 189),
    @(5188, 4946, 528,  3430, 3334, 167),
    @(4280, 3944, 622,  3396, 3247, 178),
    @(5193, 4976, 532,  6913, 6683, 284),
    @(5287, 4971, 651,  6823, 6721, 156)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j   # column B = 2 .. column G = 7
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
